# Add datasets download link
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-enter the existing F3:F12 formulas as one operation so Excel collapses
# them into a shared formula group (matches t="shared" si="0" in target).
$ws.Range("F3:F12").Formula = "=D3*E3*(E3-1)/2"

# Add the new "Download datasets" hyperlink cell in row 15.
$ws.Range("A15").Value = "Download datasets"
$ws.Hyperlinks.Add($ws.Range("A15"), "https://github.com/tridinc/genex")

# Restore the selection to match the saved view state.
$ws.Range("A20").Select()
